$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.298.80'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '3.896.66'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '524.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +7.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.11'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.611'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.82%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.718'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.72%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.172'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000332'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.89%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '41.91'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('D13').Value = '4.512.31'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.20'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.40%  '
$ws.Range('D15').Value = '3.886.56'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('E16').Value = '  +6.97%  '
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.85'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.01%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '19.66'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.79%  '
$ws.Range('D20').Value = '69.174.24'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '424.66'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('E22').Value = '  -5.70%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.15'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -6.65%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '87.70'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.27%  '
$ws.Range('E25').Value = '  +8.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.52'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.07%  '
$ws.Range('E27').Value = '  -6.58%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '36.16'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.43%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '692.96'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.47%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '13.08'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.15%  '
$ws.Range('E31').Value = '  -4.26%  '
$ws.Range('E32').Value = '  -4.11%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '68.50'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +12.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.436'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +7.94%  '
$ws.Range('E35').Value = '  -4.76%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '40.02'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.67%  '
$ws.Range('D37').Value = '0.0₃0828'
$ws.Range('E37').Value = '  -7.59%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  -3.50%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.75'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -10.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.00'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.97'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.38%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.33'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.54%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.140'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('E47').Value = '  +7.54%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0342'
$ws.Range('E48').Value = '  +2.35%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '26.65'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +4.92%  '
$ws.Range('E50').Value = '  -4.77%  '
$ws.Range('D51').Value = '2.704.28'
$ws.Range('E51').Value = '  +10.38%  '
